# Scheduled runner refresh: updates market-price-driven columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on several
# leve rows across all 8 crafter sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("I12").Value = 704.6667
$ws.Range("J12").Value = 800.75
$ws.Range("K12").Value = 704.6667
$ws.Range("L12").Value = 800.75
$ws.Range("M12").Value = -534.6667
$ws.Range("N12").Value = -1140.75
# row 15
$ws.Range("H15").Value = 831.0714
$ws.Range("I15").Value = 831.0714
$ws.Range("K15").Value = 2493.2142
$ws.Range("M15").Value = -2324.2142
# row 40
$ws.Range("H40").Value = 2969.6667
$ws.Range("I40").Value = 2910
$ws.Range("J40").Value = 2999.5
$ws.Range("K40").Value = 2910
$ws.Range("L40").Value = 2999.5
$ws.Range("M40").Value = -2735
$ws.Range("N40").Value = -3349.5
# row 105
$ws.Range("H105").Value = 77924.5
$ws.Range("J105").Value = 77924.5
$ws.Range("L105").Value = 77924.5
$ws.Range("N105").Value = -84912.5
# row 132
$ws.Range("H132").Value = 1309.8276
$ws.Range("I132").Value = 1213.75
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 3641.25
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1111.25
$ws.Range("N132").Value = -17060

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3544.66
$ws.Range("I32").Value = 2186.8108
$ws.Range("K32").Value = 2186.8108
$ws.Range("M32").Value = -1899.8108
# row 61
$ws.Range("H61").Value = 6072.375
$ws.Range("I61").Value = 3472.25
$ws.Range("J61").Value = 8672.5
$ws.Range("K61").Value = 3472.25
$ws.Range("L61").Value = 8672.5
$ws.Range("M61").Value = -3260.25
$ws.Range("N61").Value = -9096.5
# row 122
$ws.Range("H122").Value = 1846.6428
$ws.Range("I122").Value = 1839.4166
$ws.Range("K122").Value = 5518.2498
$ws.Range("M122").Value = -3068.2498
# row 132
$ws.Range("H132").Value = 2250.25
$ws.Range("I132").Value = 1404.2142
$ws.Range("K132").Value = 4212.642599999999
$ws.Range("M132").Value = -1682.642599999999
# row 136
$ws.Range("H136").Value = 6072.375
$ws.Range("I136").Value = 3472.25
$ws.Range("J136").Value = 8672.5
$ws.Range("K136").Value = 10416.75
$ws.Range("L136").Value = 26017.5
$ws.Range("M136").Value = -7866.75
$ws.Range("N136").Value = -31117.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 1440.45
$ws.Range("I20").Value = 1438.2858
$ws.Range("K20").Value = 1438.2858
$ws.Range("M20").Value = -1191.2858
# row 107
$ws.Range("H107").Value = 2927.0625
$ws.Range("I107").Value = 1969.75
$ws.Range("K107").Value = 1969.75
$ws.Range("M107").Value = -49.75
# row 134
$ws.Range("H134").Value = 6624.125
$ws.Range("I134").Value = 7035.2583
$ws.Range("K134").Value = 21105.7749
$ws.Range("M134").Value = -18570.7749

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 132
$ws.Range("H132").Value = 2297.8438
$ws.Range("I132").Value = 1183.7894
$ws.Range("J132").Value = 3926.077
$ws.Range("K132").Value = 3551.3682
$ws.Range("L132").Value = 11778.231
$ws.Range("M132").Value = -1021.3682
$ws.Range("N132").Value = -16838.231
# row 134
$ws.Range("H134").Value = 1335
$ws.Range("I134").Value = 1342.2
$ws.Range("J134").Value = 1299
$ws.Range("K134").Value = 4026.6
$ws.Range("L134").Value = 3897
$ws.Range("M134").Value = -1491.6
$ws.Range("N134").Value = -8967

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# row 33
$ws.Range("H33").Value = 82.666664
$ws.Range("I33").Value = 82.666664
$ws.Range("K33").Value = 495.999984
$ws.Range("M33").Value = -212.999984
# row 38
$ws.Range("H38").Value = 279.6
$ws.Range("I38").Value = 51
$ws.Range("K38").Value = 153
$ws.Range("M38").Value = 194
# row 107
$ws.Range("H107").Value = 650.05554
$ws.Range("J107").Value = 575.3077
$ws.Range("L107").Value = 1725.9231
$ws.Range("N107").Value = -5565.9231
# row 131
$ws.Range("H131").Value = 7824140.5
$ws.Range("J131").Value = 12793.052
$ws.Range("L131").Value = 38379.156
$ws.Range("N131").Value = -48459.156

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 146.45454
$ws.Range("I2").Value = 172.42857
$ws.Range("J2").Value = 101
$ws.Range("K2").Value = 172.42857
$ws.Range("L2").Value = 101
$ws.Range("M2").Value = -59.42857000000001
$ws.Range("N2").Value = -327
# row 97
$ws.Range("H97").Value = 931.5789
$ws.Range("J97").Value = 1972.6666
$ws.Range("L97").Value = 1972.6666
$ws.Range("N97").Value = -2964.6666
# row 102
$ws.Range("H102").Value = 2193.9
$ws.Range("I102").Value = 2482.4546
$ws.Range("J102").Value = 1841.2222
$ws.Range("K102").Value = 2482.4546
$ws.Range("L102").Value = 1841.2222
$ws.Range("M102").Value = -860.4546
$ws.Range("N102").Value = -5085.2222
# row 122
$ws.Range("H122").Value = 1513.12
$ws.Range("I122").Value = 1358.7059
$ws.Range("K122").Value = 4076.1177
$ws.Range("M122").Value = -1626.1177
# row 132
$ws.Range("H132").Value = 4573.1665
$ws.Range("I132").Value = 3161.75
$ws.Range("J132").Value = 5702.3
$ws.Range("K132").Value = 9485.25
$ws.Range("L132").Value = 17106.9
$ws.Range("M132").Value = -6955.25
$ws.Range("N132").Value = -22166.9
# row 138
$ws.Range("H138").Value = 20390
$ws.Range("I138").Value = 20390
$ws.Range("K138").Value = 20390
$ws.Range("M138").Value = -15250

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 10796.667
$ws.Range("I16").Value = 12855.8
$ws.Range("J16").Value = 501
$ws.Range("K16").Value = 12855.8
$ws.Range("L16").Value = 501
$ws.Range("M16").Value = -12685.8
$ws.Range("N16").Value = -841
# row 22
$ws.Range("H22").Value = 1008.125
$ws.Range("I22").Value = 745.8333
$ws.Range("J22").Value = 1795
$ws.Range("K22").Value = 745.8333
$ws.Range("L22").Value = 1795
$ws.Range("M22").Value = -450.8333
$ws.Range("N22").Value = -2385
# row 27
$ws.Range("H27").Value = 1008.125
$ws.Range("I27").Value = 745.8333
$ws.Range("J27").Value = 1795
$ws.Range("K27").Value = 745.8333
$ws.Range("L27").Value = 1795
$ws.Range("M27").Value = -638.8333
$ws.Range("N27").Value = -2009
# row 93
$ws.Range("H93").Value = 1449
$ws.Range("I93").Value = 1401.5
$ws.Range("J93").Value = 1496.5
$ws.Range("K93").Value = 1401.5
$ws.Range("L93").Value = 1496.5
$ws.Range("N93").Value = -3992.5
$ws.Range("M93").Value = -153.5
# row 122
$ws.Range("H122").Value = 5096.32
$ws.Range("I122").Value = 2355.3
$ws.Range("J122").Value = 6923.6665
$ws.Range("K122").Value = 7065.900000000001
$ws.Range("L122").Value = 20770.9995
$ws.Range("M122").Value = -4615.900000000001
$ws.Range("N122").Value = -25670.9995
# row 136
$ws.Range("H136").Value = 3534.1853
$ws.Range("I136").Value = 2829.6667
$ws.Range("K136").Value = 8489.000100000001
$ws.Range("M136").Value = -5939.000100000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 54870.418
$ws.Range("I122").Value = 68752.52
$ws.Range("K122").Value = 206257.56
$ws.Range("M122").Value = -203807.56
# row 132
$ws.Range("H132").Value = 3168.182
$ws.Range("I132").Value = 1340.5
$ws.Range("K132").Value = 4021.5
$ws.Range("M132").Value = -1491.5
